# Update the date line and the division practice problems in the table.
$d = $word.ActiveDocument

$replacements = @(
    @("2024-02-06 Tuesday", "2024-02-07 Wednesday"),
    @("738÷8=92, 2", "760÷6=126, 4"),
    @("879÷7=125, 4", "750÷7=107, 1"),
    @("291÷9=32, 3", "922÷7=131, 5"),
    @("111÷5=22, 1", "970÷2=485, 0"),
    @("896÷8=112, 0", "331÷9=36, 7"),
    @("605÷8=75, 5", "893÷8=111, 5"),
    @("503÷5=100, 3", "157÷8=19, 5"),
    @("408÷7=58, 2", "467÷4=116, 3"),
    @("514÷8=64, 2", "239÷8=29, 7"),
    @("362÷8=45, 2", "890÷6=148, 2"),
    @("674÷2=337, 0", "321÷2=160, 1"),
    @("314÷8=39, 2", "374÷6=62, 2"),
    @("571÷9=63, 4", "338÷3=112, 2"),
    @("524÷6=87, 2", "224÷2=112, 0"),
    @("719÷3=239, 2", "782÷7=111, 5"),
    @("142÷4=35, 2", "223÷4=55, 3"),
    @("989÷8=123, 5", "864÷7=123, 3"),
    @("658÷5=131, 3", "101÷3=33, 2"),
    @("661÷7=94, 3", "996÷3=332, 0"),
    @("329÷9=36, 5", "380÷9=42, 2"),
    @("244÷7=34, 6", "976÷5=195, 1"),
    @("743÷4=185, 3", "401÷7=57, 2"),
    @("399÷9=44, 3", "606÷5=121, 1"),
    @("383÷3=127, 2", "934÷7=133, 3"),
    @("654÷4=163, 2", "120÷3=40, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
